$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: test step now reports a timeout failure instead of a pass.
$ws.Range("L7").Value = "FAIL"
$ws.Range("M7").Value = "page.waitForTimeout: Target page, context or browser has been closed"
$ws.Range("N7").Value = "page.waitForTimeout: Target page, context or browser has been closed"
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = ""

# Rows 8-13: these steps never ran (browser/page was already closed), so
# their Status/Remarks/Actual Output/Screenshot/Page Source columns are
# cleared out entirely.
$ws.Range("L8:P13").ClearContents()
